# Apply the scheduled updates to the multiplication practice sheet.
$d = $word.ActiveDocument

# Map of old text -> new text (old values are unique within the document,
# so a plain Find/Replace for each is unambiguous).
$replacements = [ordered]@{
    "2024-02-02 Friday" = "2024-02-03 Saturday"
    "11×52=572"  = "93×11=1023"
    "94×84=7896" = "18×97=1746"
    "98×92=9016" = "21×40=840"
    "69×31=2139" = "62×73=4526"
    "69×93=6417" = "38×91=3458"
    "84×68=5712" = "20×67=1340"
    "15×76=1140" = "14×58=812"
    "37×30=1110" = "98×81=7938"
    "49×65=3185" = "41×48=1968"
    "81×50=4050" = "72×43=3096"
    "22×32=704"  = "66×16=1056"
    "97×98=9506" = "64×80=5120"
    "53×21=1113" = "96×80=7680"
    "87×68=5916" = "93×89=8277"
    "54×86=4644" = "67×49=3283"
    "48×81=3888" = "85×29=2465"
    "28×52=1456" = "24×79=1896"
    "57×65=3705" = "25×87=2175"
    "21×22=462"  = "55×41=2255"
    "61×82=5002" = "41×74=3034"
    "81×67=5427" = "68×79=5372"
    "37×61=2257" = "88×60=5280"
    "96×54=5184" = "15×35=525"
    "48×95=4560" = "83×84=6972"
    "96×98=9408" = "96×88=8448"
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host "Replaced '$old' -> '$new' : $found"
}

$d.Save()
